$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2:D6 dates (rows 2-6 get rotated values)
$ws.Range("D2").Value = 45083
$ws.Range("D3").Value = 45061
$ws.Range("D4").Value = 45084
$ws.Range("D5").Value = 45069
$ws.Range("D6").Value = 45072

# Row 2: Volumen
$ws.Range("M2").Value = 50

# Row 3: Volumen
$ws.Range("M3").Value = 40

# Row 4: Volumen / prices / precio kg
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 19000
$ws.Range("P4").Value = 18500
$ws.Range("S4").Value = 1028

# Row 5: Volumen
$ws.Range("M5").Value = 60

# Row 6: Volumen / prices / precio kg
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 15000
$ws.Range("S6").Value = 833

# Row 8: date / volumen / prices / precio kg
$ws.Range("D8").Value = 45085
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 19000
$ws.Range("O8").Value = 19000
$ws.Range("P8").Value = 19000
$ws.Range("S8").Value = 1056

# New row 9 (data previously held in row 5 before the weekly shift)
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C9").Value = "Ñuble"
$ws.Range("D9").Value = 45076
$ws.Range("D9").NumberFormat = $ws.Range("D8").NumberFormat
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100104
$ws.Range("H9").Value = "Frutos de pepita"
$ws.Range("I9").Value = 100104001
$ws.Range("J9").Value = "Granada"
$ws.Range("K9").Value = "Wonderfull"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 20
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 15000
$ws.Range("Q9").Value = "$/caja 18 kilos granel"
$ws.Range("R9").Value = "Provincia de Curicó"
$ws.Range("S9").Value = 833
$ws.Range("T9").Value = 18
